$d = $word.ActiveDocument

# The commit reorders/extends the bullet list. Easiest + most faithful way
# to reproduce the exact post-edit run/proofErr/bookmark structure is to
# rebuild the whole body from scratch as WordprocessingML and drop it in,
# rather than trying to replicate every paragraph move with Find/Replace.

# Collapse the document down to a single empty paragraph.
while ($d.Paragraphs.Count -gt 1) {
    $d.Paragraphs($d.Paragraphs.Count).Range.Delete()
}
$d.Range(0, $d.Content.End).Text = ""

# Each entry is one target <w:p> (post-edit order). Spelling-flagged
# tokens keep their <w:proofErr> wrap; plain bullets are single runs.
$paragraphs = @(
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:r><w:t>Controller Service Repository Layers</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Jsp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as view</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Use</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t xml:space="preserve"> Validator</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:r><w:t>Errors</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used message resources for Labels</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used message resources</w:t></w:r><w:r><w:t xml:space="preserve"> for Errors</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CrudRepository</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
    '<w:p><w:r><w:t>Used JPA</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Enum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for Gender</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Comparator.comparing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for double time sorting for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>praentId</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>serviced</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used Spring Security</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>db</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Authentication</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used Role based Authentication</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used Role based home page</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used High charts</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used Ajax</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used Lightbox popup</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used native </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> queries</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used Global Exception Handler</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used Custom Error Control</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mysql</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Logback</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
    '<w:p><w:r><w:t>Used bootstrap</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used DTO pattern</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used Converter class from Entity to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Viceversa</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
    '<w:p><w:r><w:t>Used Tiles</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used Interceptor</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used filter</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used custom Exception types</w:t></w:r></w:p>',
    '<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Used Transaction Management </w:t></w:r><w:r><w:t>(@Transactional)</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Removed </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>JSessionId</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">Used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Bcypt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Password Encoder</w:t></w:r></w:p>',
    '<w:p><w:r><w:t>Used Profiles</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
)

$bodyXml = $paragraphs -join ""
$pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + 
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + 
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + 
    $bodyXml + 
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r = $d.Range(0, $d.Content.End)
$r.InsertXML($pkgXml)

